$wb = $excel.ActiveWorkbook

# =====================================================================
# Step 1: "总计" (Total) sheet - insert a new row 2 for the "2022-Q4"
# quarter, shifting the existing rows (old 2022-Q3 .. 2020-Q4) down by
# one.
# =====================================================================
$totalSheet = $wb.Worksheets.Item(1)
$totalSheet.Rows.Item(2).Insert()

# Row 3 (the old row 2, shifted down by the insert above) already has
# the correct per-column styling for this table; clone its column-A
# cell (bordered/bold index style) onto the new row 2, then clear the
# incidental style that Insert() propagated onto B2:D2 from the header
# row above, so the new row matches its siblings exactly.
$totalSheet.Range("A3").Copy($totalSheet.Range("A2"))
$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 39
$totalSheet.Range("D2").Value = 3.19

# =====================================================================
# Step 2: insert a brand-new "2022-Q4" worksheet right after "总计" and
# before the existing "2022-Q3" sheet, carrying the per-fund breakdown.
# =====================================================================
$q4Sheet = $wb.Worksheets.Add($null, $totalSheet)
$q4Sheet.Name = "2022-Q4"

# The sibling quarter sheets all share identical header-row styling
# (bold font + border on B1:H1) and an identical bordered/bold style on
# the column-A index. Re-fetch "2022-Q3" now that it has been pushed to
# position 3, and clone both styles from it onto the new sheet before
# filling in values, so no ad-hoc style entries get invented.
$q3Sheet = $wb.Worksheets.Item(3)
$q3Sheet.Range("B1:H1").Copy($q4Sheet.Range("B1:H1"))

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q4Sheet.Cells.Item(1, 2 + $i).Value = $headers[$i]
}

$data2022Q4 = @(
    ,('0','009861','鹏华新兴成长混合A','19.32','94.27','3.38','0.6530','6')
    ,('1','005454','前海开源医疗健康灵活配置混合C','11.09','87.03','4.56','0.5057','7')
    ,('2','159883','永赢中证全指医疗器械ETF','22.83','99.25','2.13','0.4863','10')
    ,('3','005453','前海开源医疗健康灵活配置混合A','9.08','87.03','4.56','0.4140','7')
    ,('4','005028','鹏华研究精选灵活配置混合','11.55','94.15','3.39','0.3915','8')
    ,('5','013967','广发恒享一年持有期混合A','19.39','20.90','0.52','0.1008','5')
    ,('6','011876','景顺长城医疗健康混合A','2.34','86.33','4.28','0.1002','6')
    ,('7','011601','前海开源公共卫生主题精选股票A','1.26','85.69','4.52','0.0570','7')
    ,('8','011877','景顺长城医疗健康混合C','1.23','86.33','4.28','0.0526','6')
    ,('9','013968','广发恒享一年持有期混合C','9.93','20.90','0.52','0.0516','5')
    ,('10','010029','富国稳进回报12个月持有期混合A','6.03','20.60','0.80','0.0482','9')
    ,('11','002152','华宝核心优势灵活配置混合A','0.44','92.49','9.37','0.0412','1')
    ,('12','011602','前海开源公共卫生主题精选股票C','0.82','85.69','4.52','0.0371','7')
    ,('13','159898','招商中证全指医疗器械ETF','1.70','99.16','2.13','0.0362','10')
    ,('14','009862','鹏华新兴成长混合C','0.64','94.27','3.38','0.0216','6')
    ,('15','164401','前海开源中证健康产业指数','1.89','94.43','1.14','0.0215','3')
    ,('16','007424','西部利得聚禾灵活配置混合C','0.56','80.30','3.81','0.0213','8')
    ,('17','005360','汇安资产轮动灵活配置混合A','0.26','94.03','6.63','0.0172','9')
    ,('18','014206','长城中证医药卫生指数增强C','0.38','94.69','3.75','0.0142','7')
    ,('19','002271','招商安弘灵活配置混合','0.55','72.78','2.27','0.0125','8')
    ,('20','159797','汇添富中证全指医疗器械ETF','0.56','98.45','2.13','0.0119','10')
    ,('21','014205','长城中证医药卫生指数增强A','0.27','94.69','3.75','0.0101','7')
    ,('22','003659','山西证券策略精选灵活配置混合','0.50','50.39','1.80','0.0090','3')
    ,('23','015507','兴业中证500指数增强A','0.71','92.30','1.22','0.0087','8')
    ,('24','010154','中加中证500指数增强C','0.51','94.31','1.65','0.0084','5')
    ,('25','002135','广发鑫源灵活配置混合A','0.59','27.74','1.41','0.0083','6')
    ,('26','010153','中加中证500指数增强A','0.46','94.31','1.65','0.0076','5')
    ,('27','010030','富国稳进回报12个月持有期混合C','0.89','20.60','0.80','0.0071','9')
    ,('28','006818','安信盈利驱动股票A','0.08','88.09','6.25','0.0050','1')
    ,('29','510560','国寿安保中证500ETF','1.81','99.21','0.25','0.0045','8')
    ,('30','015508','兴业中证500指数增强C','0.37','92.30','1.22','0.0045','8')
    ,('31','006819','安信盈利驱动股票C','0.06','88.09','6.25','0.0038','1')
    ,('32','851088','海通量化成长精选一年持有期混合A','0.39','83.15','0.93','0.0036','4')
    ,('33','007423','西部利得聚禾灵活配置混合A','0.08','80.30','3.81','0.0030','8')
    ,('34','002136','广发鑫源灵活配置混合C','0.19','27.74','1.41','0.0027','6')
    ,('35','850010','海通量化成长精选一年持有期混合B','0.26','83.15','0.93','0.0024','4')
    ,('36','016461','华宝核心优势灵活配置混合C','0.02','92.49','9.37','0.0019','1')
    ,('37','017213','汇安资产轮动灵活配置混合C','0.01','94.03','6.63','0.0007','9')
    ,('38','851099','海通量化成长精选一年持有期混合C','0.03','83.15','0.93','0.0003','4')
)

for ($i = 0; $i -lt $data2022Q4.Length; $i++) {
    $r = $i + 2
    $rec = $data2022Q4[$i]

    # Column A: numeric 0-based index, styled like the rest of the
    # table (clone from 2022-Q3's own index column, then overwrite the
    # value with the correct index).
    $q3Sheet.Range("A2").Copy($q4Sheet.Cells.Item($r, 1))
    $q4Sheet.Cells.Item($r, 1).Value = [double]$rec[0]

    # Columns B-G: text, even when the literal looks numeric (fund
    # codes with leading zeros, percentages with trailing zeros, etc).
    # Force text storage via NumberFormat "@" before assigning, then
    # clear the incidental format so no stray style survives the round
    # trip.
    for ($c = 1; $c -le 6; $c++) {
        $cell = $q4Sheet.Cells.Item($r, 1 + $c)
        $cell.NumberFormat = "@"
        $cell.Value = [string]$rec[$c]
        $cell.ClearFormats()
    }

    # Column H: numeric rank.
    $q4Sheet.Cells.Item($r, 8).Value = [double]$rec[7]
}

Write-Host "done"
